$wb = $excel.ActiveWorkbook

$sheet1Urls = @(
    "https://ainow.ai/2019/11/25/181157/",
    "https://ainow.ai/2019/11/20/181109/",
    "https://ainow.ai/2019/11/20/181105/",
    "https://ainow.ai/2019/11/20/181096/",
    "https://ainow.ai/2019/11/20/180352/",
    "https://ainow.ai/2019/11/19/180997/",
    "https://ainow.ai/2019/11/19/180379/",
    "https://ainow.ai/2019/11/18/180975/",
    "https://ainow.ai/2019/11/18/180958/",
    "https://ainow.ai/2019/11/18/180961/",
    "https://ainow.ai/2019/11/18/180936/",
    "https://ainow.ai/2019/11/18/180949/",
    "https://ainow.ai/2019/11/14/180803/",
    "https://ainow.ai/2019/11/14/180805/",
    "https://ainow.ai/2019/11/14/180757/",
    "https://ainow.ai/2019/11/14/180721/",
    "https://ainow.ai/2019/11/13/180796/",
    "https://ainow.ai/2019/11/13/180323/",
    "https://ainow.ai/2019/11/12/180479/",
    "https://ainow.ai/2019/11/09/180524/"
)

$sheet2Urls = @(
    "https://ainow.ai/2019/09/30/175475/",
    "https://ainow.ai/2019/06/13/172100/",
    "https://ainow.ai/2019/06/09/171839/",
    "https://ainow.ai/2019/06/07/171703/",
    "https://ainow.ai/2019/06/04/171382/",
    "https://ainow.ai/2019/05/30/170999/",
    "https://ainow.ai/2019/05/29/170919/",
    "https://ainow.ai/2019/05/29/170915/",
    "https://ainow.ai/2019/05/29/170910/",
    "https://ainow.ai/2019/05/24/170530/",
    "https://ainow.ai/2019/05/22/170362/",
    "https://ainow.ai/2019/05/20/170150/",
    "https://ainow.ai/2019/05/20/170141/",
    "https://ainow.ai/2019/05/16/169893/",
    "https://ainow.ai/2019/05/15/169799/",
    "https://ainow.ai/2019/05/15/169797/",
    "https://ainow.ai/2019/05/15/169796/",
    "https://ainow.ai/2019/05/12/169562/",
    "https://ainow.ai/2019/05/08/169194/",
    "https://ainow.ai/2019/04/27/168687/"
)

$wsAI = $wb.Worksheets.Item("AI")
$row = 23
foreach ($url in $sheet1Urls) {
    $cell = $wsAI.Cells.Item($row, 1)
    $cell.Value = $url
    $cell.Style = "Normal"
    $row++
}

$wsIoT = $wb.Worksheets.Item("IoT")
$row = 23
foreach ($url in $sheet2Urls) {
    $cell = $wsIoT.Cells.Item($row, 1)
    $cell.Value = $url
    $cell.Style = "Normal"
    $row++
}
